$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert 3 new rows right before the old last data row (row 22),
#    pushing it down to row 25 and the footer block (27/28) to 30/31.
# ------------------------------------------------------------------
$ws.Rows("22:24").Insert()

# ------------------------------------------------------------------
# 2) Copy the formatting of a "normal" data row onto the new rows
#    so borders / fonts / number formats match the rest of the table.
# ------------------------------------------------------------------
$ws.Range("B16:J16").Copy()
$ws.Range("B22:J24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) New period "2508" data: re-use the same five workers that were
#    billed for period 2507, shifting the old row 21/22 duplicate
#    (GLENDA MARCELA MIELES GOMEZ / period 2301 & 2212) out.
# ------------------------------------------------------------------
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "73006146"
$ws.Range("D21").Value = "CEDRICK CONTRERA GUARDO"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 128000
$ws.Range("G21").Value = 3200000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "45687489"
$ws.Range("D22").Value = "XIOMARA PATRICIA MONROY TINOCO"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 128000
$ws.Range("G22").Value = 3200000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "32939066"
$ws.Range("D23").Value = "KAREN MARGARITA SARABIA AYOLA"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1047475016"
$ws.Range("D24").Value = "NELSON ENRIQUE GOMEZ VEGA"
$ws.Range("E24").Value = "2508"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1047434781"
$ws.Range("D25").Value = "DARWIN ALEXANDER ESCOBAR MIRANDA"
$ws.Range("E25").Value = "2508"
$ws.Range("F25").Value = 56940
$ws.Range("G25").Value = 1423500

# ------------------------------------------------------------------
# 4) Refresh the summary figures at the top of the sheet.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 853640
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 2
